$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.443521022796631
$ws.Range("B1").Value = 4.901392936706543
$ws.Range("C1").Value = 7.932816505432129
$ws.Range("D1").Value = 7.517411231994629
$ws.Range("E1").Value = 4.420957088470459
